# Apply "update course list and research experience" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 11: fill in a new "Research" entry (SRP199, 01/2022, grade A+, 2 credits) ---
$ws.Range("H11").Value = "SRP199"
$ws.Range("J11").Value = "01/2022"
$ws.Range("L11").Value = "A+"
$ws.Range("M11").Value = 2

# --- "Currently Taking" block (rows 27-31): drop the header row, shift the four
#     course rows up by one, add completed grades, and clear the now-unused row 31 ---
$ws.Range("A27").Value = "CS118 "
$ws.Range("B27").Value = "Computer Network Fundamentals"
$ws.Range("C27").Value = "01/2022"
$ws.Range("D27").Value = "CS"
$ws.Range("E27").Value = "A+"
$ws.Range("F27").Value = 4
$ws.Range("A27").Font.Bold = $false
# B27:F27 are brand-new cells - match the regular data-row font (style index 1)
# instead of leaving them on the workbook default style.
$ws.Range("B27:F27").Font.Name = "Times Roman"
$ws.Range("B27:F27").Font.Size = 12

$ws.Range("A28").Value = "CS188"
$ws.Range("B28").Value = "Natural Language Processing"
$ws.Range("C28").Value = "01/2022"
$ws.Range("D28").Value = "CS"
$ws.Range("E28").Value = "A"
$ws.Range("F28").Value = 4
$ws.Range("E28").Font.Name = "Times Roman"
$ws.Range("E28").Font.Size = 12

$ws.Range("A29").Value = "CS188"
$ws.Range("B29").Value = "Deep Learning In Computer Vision"
$ws.Range("C29").Value = "01/2022"
$ws.Range("D29").Value = "CS"
$ws.Range("E29").Value = "A"
$ws.Range("F29").Value = 4

$ws.Range("A30").Value = "CSM148"
$ws.Range("B30").Value = "Data Science"
$ws.Range("C30").Value = "01/2022"
$ws.Range("D30").Value = "CS"
$ws.Range("E30").Value = "A+"
$ws.Range("F30").Value = 4

$ws.Range("A31").Value = ""
$ws.Range("B31").Value = ""
$ws.Range("C31").Value = ""
$ws.Range("D31").Value = ""
$ws.Range("E31").Value = ""
$ws.Range("F31").Value = ""

# --- "Planning to take" section becomes "Currently taking" (research experience update) ---
$ws.Range("A33").Value = "Currently taking"

# Replace the old CSM146 row with the new ECE188 course
$ws.Range("A35").Value = "ECE188"
$ws.Range("B35").Value = "Secure Computing Systems"
$ws.Range("D35").Value = "ECE"

# CS130 / Software Engineering row now has a start date and credit value
$ws.Range("C34").Value = "03/2022"
$ws.Range("F34").Value = 4
$ws.Range("C35").Value = "03/2022"

# Remove the CS133 and CS136 rows entirely (no longer planned)
$ws.Range("A36").Value = ""
$ws.Range("B36").Value = ""
$ws.Range("C36").Value = ""
$ws.Range("D36").Value = ""

$ws.Range("A37").Value = ""
$ws.Range("B37").Value = ""
$ws.Range("C37").Value = ""
$ws.Range("D37").Value = ""

# --- Restore the active-cell selection to match where the author left off editing ---
$ws.Range("I28").Select()
